$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset's metadata for the "provincia" column is re-curated:
# - D2: sdmx-dimension:refArea      -> iaest-measure:provincia
# - D3: dim                         -> medida
# - D4: URI-Provincia               -> xsd:int
$ws.Range("D2").Value = "iaest-measure:provincia"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"
